{"js": "// Office.js (Word JavaScript API) script\n// Applies the LOM3256 content update described by the diff:\n//  - shorten the course title / subtitle (drop \"dos Materiais\" / \"of materials\")\n//  - update the \"Ativa\u00e7\u00e3o\" date\n//  - collapse the line-broken \"Objetivos\" paragraph into one run and\n//    append its English translation as a new italic paragraph\n//  - collapse the \"Programa resumido\" paragraph into one run and\n//    append its English translation as a new italic paragraph\n//  - collapse the long \"Programa\" bullet paragraph into one run and\n//    append its English translation as a new italic paragraph\n//  - update the \"Requisitos\" bullet entry\n\nconst body = context.document.body;\n\n// 1) Title heading.\nconst titleResults = body.search(\n  \"LOM3256 -  T\u00f3picos em C\u00e1lculo de Estrutura Eletr\u00f4nica dos Materiais\",\n  { matchCase: true }\n);\ntitleResults.load(\"items\");\nawait context.sync();\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\n    \"LOM3256 -  T\u00f3picos em C\u00e1lculo de Estrutura Eletr\u00f4nica\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 2) English subtitle heading.\nconst subtitleResults = body.search(\n  \"Methods of electronic structure calculation of materials\",\n  { matchCase: true }\n);\nsubtitleResults.load(\"items\");\nawait context.sync();\nif (subtitleResults.items.length > 0) {\n  subtitleResults.items[0].insertText(\n    \"Methods of electronic structure calculation\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 3) Ativa\u00e7\u00e3o date.\nconst dateResults = body.search(\"Ativa\u00e7\u00e3o: 15/07/2015\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\n    \"Ativa\u00e7\u00e3o: 01/01/2023\",\n    Word.InsertLocation.replace\n  );\n}\n\n// Helper: find the paragraph that immediately follows a heading paragraph\n// whose text exactly equals `headingText`. Re-loads the paragraph\n// collection fresh so indices stay valid across edits.\nasync function paragraphAfterHeading(headingText) {\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n  const idx = paragraphs.items.findIndex((p) => p.text === headingText);\n  if (idx === -1 || idx + 1 >= paragraphs.items.length) {\n    return null;\n  }\n  return paragraphs.items[idx + 1];\n}\n\n// 4) Objetivos: merge the broken-up runs into a single paragraph of text,\n//    then add the English translation as a new italic paragraph after it.\nconst objetivosPt =\n  \"Propiciar ao aluno uma vis\u00e3o b\u00e1sica sobre os principais m\u00e9todos de determina\u00e7\u00e3o te\u00f3rica da estrutura eletr\u00f4nica dos materiais, com enfoque em s\u00f3lidos cristalinos, mas tamb\u00e9m em materiais bidimensionais e nanoestruturados.O principal m\u00e9todo de c\u00e1lculo a ser empregado no curso ser\u00e1 a Teoria do Funcional da Densidade(Density Functional Theory, DFT), em algumas de suas muitas variantes. Ao final do curso, o aluno estar\u00e1 apto a determinar propriedades dos materiais como estruturas de bandas, densidades de estados, superf\u00edcies de Fermi e constantes el\u00e1sticas, usando um ou mais dos m\u00e9todos e c\u00f3digos computacionais apresentados em aula.\";\nconst objetivosEn =\n  \"Provide the student with a basic view of the main methods of theoretical determination of the electronic structure, focusing on crystalline solids, but also on molecules, two-dimensional materials and nanostructured materials. The main calculation method to be used in the course will be the Density Functional Theory (DFT), in some of its many variants. At the end of the course, the student will be able to determine material properties such as band structures, densities of states, elastic constants, and Fermi surfaces, using one or more of the methods and computer codes presented in class.\";\n\nlet targetPara = await paragraphAfterHeading(\"Objetivos\");\nif (targetPara) {\n  targetPara.clear();\n  targetPara.insertText(objetivosPt, Word.InsertLocation.start);\n  const newPara = targetPara.insertParagraph(\n    objetivosEn,\n    Word.InsertLocation.after\n  );\n  newPara.font.italic = true;\n  await context.sync();\n}\n\n// 5) Programa resumido: merge the two runs, then add the English translation.\nconst resumidoPt =\n  \"Revis\u00e3o de mec\u00e2nica qu\u00e2ntica; Revis\u00e3o de f\u00edsica do estado s\u00f3lido; M\u00e9todo de Hartree-Fock; Teoria do funcional da densidade; M\u00e9todos de ondas planas e pseudo-potenciais; C\u00f3digos computacionais\";\nconst resumidoEn =\n  \"Review of Quantum Mechanics; Review of Solid State Physics; Hartree-Fock Method; Density Functional Theory; Plane and pseudopotential wave methods; computer codes\";\n\ntargetPara = await paragraphAfterHeading(\"Programa resumido\");\nif (targetPara) {\n  targetPara.clear();\n  targetPara.insertText(resumidoPt, Word.InsertLocation.start);\n  const newPara = targetPara.insertParagraph(\n    resumidoEn,\n    Word.InsertLocation.after\n  );\n  newPara.font.italic = true;\n  await context.sync();\n}\n\n// 6) Programa: merge the long bullet list into one run, then add the\n//    English translation.\nconst programaPt =\n  \"Revis\u00e3o de mec\u00e2nica qu\u00e2nticao Equa\u00e7\u00e3o de Schr\u00f6dingero \u00c1tomo do hidrog\u00eanio e orbitais at\u00f4micoso Nota\u00e7\u00e3o de Diraco Princ\u00edpio variacionalo Combina\u00e7\u00e3o linear de orbitais at\u00f4micosRevis\u00e3o de f\u00edsica do estado s\u00f3lidoo Espa\u00e7o direto e rec\u00edprocoo Teorema de Blocho Zona de Brillouino Bandas de energia e densidade de estadoso Energia de Fermi e superf\u00edcie de Fermio Aproxima\u00e7\u00e3o de el\u00e9trons livresM\u00e9todo de Hartree-Focko Determinantes de Slatero Equa\u00e7\u00e3o de Hartree-Focko Potencial de troca e correla\u00e7\u00e3oo Algoritmo autoconsistenteTeoria do funcional da densidadeo Teoremas de Hohenberg-Kohno Equa\u00e7\u00f5es de Kohn-Shamo Funcionais de troca e correla\u00e7\u00e3o: LDA, GGA, etc.M\u00e9todos de ondas planas e pseudo-potenciaiso Bases de ondas planaso Pseudo-potenciaiso Bases de ondas planas aumentadas e linearizadaso M\u00e9todo FP-LAPWC\u00f3digos computacionaiso Quantum Espressoo Elko Wien2ko VASP\";\nconst programaEn =\n  \"\u2022 Review of quantum mechanics: Schr\u00f6dinger's equation; Hydrogen atom and atomic orbitals; Dirac notation; Variational principle; Linear combination of atomic orbitals. \u2022 Solid state physics review: Direct and reciprocal space; Bloch's Theorem; Brillouin zone; Energy bands and density of states; Fermi energy and Fermi surface; Free electrons Approximation. \u2022 Hartree-Fock method: Slater determinants; Hartree-Fock equation; Exchange and correlation potential; Self-consistent algorithm. \u2022 Density functional theory: Hohenberg-Kohn theorems; Kohn-Sham equations; Exchange and correlation functionals: LDA, GGA, etc. \u2022 Plane and pseudopotential wave methods: Plane wave bases; Pseudo-potentials; \u2022 Augmented and linearized plane wave bases: FP-LAPW method. \u2022 Computer codes: NWCHEM, Quantum Espresso, , Wien2k, exciting, VASP, etc.\";\n\ntargetPara = await paragraphAfterHeading(\"Programa\");\nif (targetPara) {\n  targetPara.clear();\n  targetPara.insertText(programaPt, Word.InsertLocation.start);\n  const newPara = targetPara.insertParagraph(\n    programaEn,\n    Word.InsertLocation.after\n  );\n  newPara.font.italic = true;\n  await context.sync();\n}\n\n// 7) Requisitos bullet entry.\nconst reqResults = body.search(\n  \"LOM3215 -  F\u00edsica do Estado S\u00f3lido  (Requisito)\",\n  { matchCase: true }\n);\nreqResults.load(\"items\");\nawait context.sync();\nif (reqResults.items.length > 0) {\n  reqResults.items[0].insertText(\n    \"LOM3226 -  Mec\u00e2nica Qu\u00e2ntica  (Requisito)\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Applies the LOM3256 content update described by the diff:\n#  - shorten the course title / subtitle (drop \"dos Materiais\" / \"of materials\")\n#  - update the \"Ativa\u00e7\u00e3o\" date\n#  - collapse the line-broken \"Objetivos\" paragraph into one run and\n#    append its English translation as a new italic paragraph\n#  - collapse the \"Programa resumido\" paragraph into one run and\n#    append its English translation as a new italic paragraph\n#  - collapse the long \"Programa\" bullet paragraph into one run and\n#    append its English translation as a new italic paragraph\n#  - update the \"Requisitos\" bullet entry\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndexByExactText($doc, $text) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $t = $doc.Paragraphs.Item($i).Range.Text\n        $t = $t.TrimEnd([char]13, [char]7)\n        if ($t -eq $text) { return $i }\n    }\n    return -1\n}\n\n# 1) Title heading.\n$find = $d.Content.Find\n$find.Execute(\n    \"LOM3256 -  T\u00f3picos em C\u00e1lculo de Estrutura Eletr\u00f4nica dos Materiais\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"LOM3256 -  T\u00f3picos em C\u00e1lculo de Estrutura Eletr\u00f4nica\", 2\n) | Out-Null\n\n# 2) English subtitle heading.\n$find = $d.Content.Find\n$find.Execute(\n    \"Methods of electronic structure calculation of materials\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"Methods of electronic structure calculation\", 2\n) | Out-Null\n\n# 3) Ativa\u00e7\u00e3o date.\n$find = $d.Content.Find\n$find.Execute(\n    \"Ativa\u00e7\u00e3o: 15/07/2015\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"Ativa\u00e7\u00e3o: 01/01/2023\", 2\n) | Out-Null\n\n# 4) Objetivos: merge the broken-up runs into a single paragraph of text,\n#    then add the English translation as a new italic paragraph after it.\n$objetivosPt = \"Propiciar ao aluno uma vis\u00e3o b\u00e1sica sobre os principais m\u00e9todos de determina\u00e7\u00e3o te\u00f3rica da estrutura eletr\u00f4nica dos materiais, com enfoque em s\u00f3lidos cristalinos, mas tamb\u00e9m em materiais bidimensionais e nanoestruturados.O principal m\u00e9todo de c\u00e1lculo a ser empregado no curso ser\u00e1 a Teoria do Funcional da Densidade(Density Functional Theory, DFT), em algumas de suas muitas variantes. Ao final do curso, o aluno estar\u00e1 apto a determinar propriedades dos materiais como estruturas de bandas, densidades de estados, superf\u00edcies de Fermi e constantes el\u00e1sticas, usando um ou mais dos m\u00e9todos e c\u00f3digos computacionais apresentados em aula.\"\n$objetivosEn = \"Provide the student with a basic view of the main methods of theoretical determination of the electronic structure, focusing on crystalline solids, but also on molecules, two-dimensional materials and nanostructured materials. The main calculation method to be used in the course will be the Density Functional Theory (DFT), in some of its many variants. At the end of the course, the student will be able to determine material properties such as band structures, densities of states, elastic constants, and Fermi surfaces, using one or more of the methods and computer codes presented in class.\"\n\n$idx = Find-ParagraphIndexByExactText $d \"Objetivos\"\nif ($idx -gt 0) {\n    $contentPara = $d.Paragraphs.Item($idx + 1)\n    $contentPara.Range.Text = $objetivosPt\n    $contentPara.Range.InsertParagraphAfter()\n    $newPara = $contentPara.Next()\n    $newRange = $newPara.Range\n    $newRange.InsertAfter($objetivosEn)\n    $newRange.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1; exclude the paragraph mark\n    $newRange.Font.Italic = $true\n}\n\n# 5) Programa resumido: merge the two runs, then add the English translation.\n$resumidoPt = \"Revis\u00e3o de mec\u00e2nica qu\u00e2ntica; Revis\u00e3o de f\u00edsica do estado s\u00f3lido; M\u00e9todo de Hartree-Fock; Teoria do funcional da densidade; M\u00e9todos de ondas planas e pseudo-potenciais; C\u00f3digos computacionais\"\n$resumidoEn = \"Review of Quantum Mechanics; Review of Solid State Physics; Hartree-Fock Method; Density Functional Theory; Plane and pseudopotential wave methods; computer codes\"\n\n$idx = Find-ParagraphIndexByExactText $d \"Programa resumido\"\nif ($idx -gt 0) {\n    $contentPara = $d.Paragraphs.Item($idx + 1)\n    $contentPara.Range.Text = $resumidoPt\n    $contentPara.Range.InsertParagraphAfter()\n    $newPara = $contentPara.Next()\n    $newRange = $newPara.Range\n    $newRange.InsertAfter($resumidoEn)\n    $newRange.MoveEnd(1, -1) | Out-Null\n    $newRange.Font.Italic = $true\n}\n\n# 6) Programa: merge the long bullet list into one run, then add the\n#    English translation.\n$programaPt = \"Revis\u00e3o de mec\u00e2nica qu\u00e2nticao Equa\u00e7\u00e3o de Schr\u00f6dingero \u00c1tomo do hidrog\u00eanio e orbitais at\u00f4micoso Nota\u00e7\u00e3o de Diraco Princ\u00edpio variacionalo Combina\u00e7\u00e3o linear de orbitais at\u00f4micosRevis\u00e3o de f\u00edsica do estado s\u00f3lidoo Espa\u00e7o direto e rec\u00edprocoo Teorema de Blocho Zona de Brillouino Bandas de energia e densidade de estadoso Energia de Fermi e superf\u00edcie de Fermio Aproxima\u00e7\u00e3o de el\u00e9trons livresM\u00e9todo de Hartree-Focko Determinantes de Slatero Equa\u00e7\u00e3o de Hartree-Focko Potencial de troca e correla\u00e7\u00e3oo Algoritmo autoconsistenteTeoria do funcional da densidadeo Teoremas de Hohenberg-Kohno Equa\u00e7\u00f5es de Kohn-Shamo Funcionais de troca e correla\u00e7\u00e3o: LDA, GGA, etc.M\u00e9todos de ondas planas e pseudo-potenciaiso Bases de ondas planaso Pseudo-potenciaiso Bases de ondas planas aumentadas e linearizadaso M\u00e9todo FP-LAPWC\u00f3digos computacionaiso Quantum Espressoo Elko Wien2ko VASP\"\n$programaEn = \"\u2022 Review of quantum mechanics: Schr\u00f6dinger's equation; Hydrogen atom and atomic orbitals; Dirac notation; Variational principle; Linear combination of atomic orbitals. \u2022 Solid state physics review: Direct and reciprocal space; Bloch's Theorem; Brillouin zone; Energy bands and density of states; Fermi energy and Fermi surface; Free electrons Approximation. \u2022 Hartree-Fock method: Slater determinants; Hartree-Fock equation; Exchange and correlation potential; Self-consistent algorithm. \u2022 Density functional theory: Hohenberg-Kohn theorems; Kohn-Sham equations; Exchange and correlation functionals: LDA, GGA, etc. \u2022 Plane and pseudopotential wave methods: Plane wave bases; Pseudo-potentials; \u2022 Augmented and linearized plane wave bases: FP-LAPW method. \u2022 Computer codes: NWCHEM, Quantum Espresso, , Wien2k, exciting, VASP, etc.\"\n\n$idx = Find-ParagraphIndexByExactText $d \"Programa\"\nif ($idx -gt 0) {\n    $contentPara = $d.Paragraphs.Item($idx + 1)\n    $contentPara.Range.Text = $programaPt\n    $contentPara.Range.InsertParagraphAfter()\n    $newPara = $contentPara.Next()\n    $newRange = $newPara.Range\n    $newRange.InsertAfter($programaEn)\n    $newRange.MoveEnd(1, -1) | Out-Null\n    $newRange.Font.Italic = $true\n}\n\n# 7) Requisitos bullet entry.\n$find = $d.Content.Find\n$find.Execute(\n    \"LOM3215 -  F\u00edsica do Estado S\u00f3lido  (Requisito)\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"LOM3226 -  Mec\u00e2nica Qu\u00e2ntica  (Requisito)\", 2\n) | Out-Null\n"}
